$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 2882
$ws.Range("D2").Value = 2842
$ws.Range("E2").Value = 1487
$ws.Range("F2").Value = 2882
$ws.Range("G2").Value = 2431
$ws.Range("H2").Value = 2251
$ws.Range("I2").Value = 2250
$ws.Range("J2").Value = 2251
$ws.Range("K2").Value = 2430
$ws.Range("L2").Value = 2882
$ws.Range("G3").Value = 190
$ws.Range("H3").Value = 190
$ws.Range("I3").Value = 190
$ws.Range("J3").Value = 190
$ws.Range("K3").Value = 190
$ws.Range("C4").Value = 748
$ws.Range("E4").Value = 672
$ws.Range("F4").Value = 747
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 628
$ws.Range("I4").Value = 628
$ws.Range("J4").Value = 628
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 747
$ws.Range("G5").Value = 465
$ws.Range("H5").Value = 465
$ws.Range("I5").Value = 465
$ws.Range("J5").Value = 465
$ws.Range("K5").Value = 465
$ws.Range("L5").Value = 395
$ws.Range("B6").Value = 1038
$ws.Range("C6").Value = 520
$ws.Range("D6").Value = 1038
$ws.Range("E6").Value = 520
$ws.Range("G6").Value = 767
$ws.Range("H6").Value = 800
$ws.Range("I6").Value = 751
$ws.Range("J6").Value = 751
$ws.Range("K6").Value = 721
$ws.Range("L6").Value = 520
$ws.Range("B7").Value = 1365
$ws.Range("C7").Value = 1363
$ws.Range("D7").Value = 1365
$ws.Range("E7").Value = 1360
$ws.Range("F7").Value = 1452
$ws.Range("H7").Value = 1452
$ws.Range("I7").Value = 1365
$ws.Range("J7").Value = 1365
$ws.Range("K7").Value = 1365
$ws.Range("L7").Value = 1452
$ws.Range("B8").Value = 511
$ws.Range("C8").Value = 511
$ws.Range("D8").Value = 511
$ws.Range("E8").Value = 510
$ws.Range("F8").Value = 554
$ws.Range("G8").Value = 554
$ws.Range("I8").Value = 511
$ws.Range("J8").Value = 511
$ws.Range("K8").Value = 511
$ws.Range("L8").Value = 510
$ws.Range("G9").Value = 408
$ws.Range("K9").Value = 408
$ws.Range("G10").Value = 310
$ws.Range("H11").Value = 443
$ws.Range("I11").Value = 443
$ws.Range("G12").Value = 767
$ws.Range("H12").Value = 861
$ws.Range("I12").Value = 861
$ws.Range("J12").Value = 861
$ws.Range("K12").Value = 767
